# Scheduled runner update: refresh market-board price/profit figures
# across the per-job Leve profit sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 33074
$ws.Range("J47").Value = 33074
$ws.Range("L47").Value = 33074
$ws.Range("N47").Value = -35018

$ws.Range("H64").Value = 2931.8
$ws.Range("I64").Value = 2900
$ws.Range("K64").Value = 2900
$ws.Range("M64").Value = -2652

$ws.Range("H67").Value = 2931.8
$ws.Range("I67").Value = 2900
$ws.Range("K67").Value = 2900
$ws.Range("M67").Value = -2042

$ws.Range("H132").Value = 916505.4399999999
$ws.Range("I132").Value = 979115.75
$ws.Range("J132").Value = 2395
$ws.Range("K132").Value = 2937347.25
$ws.Range("L132").Value = 7185
$ws.Range("M132").Value = -2934817.25
$ws.Range("N132").Value = -12245

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2983.077
$ws.Range("I88").Value = 1490
$ws.Range("J88").Value = 3107.5
$ws.Range("K88").Value = 1490
$ws.Range("L88").Value = 3107.5
$ws.Range("M88").Value = -1084
$ws.Range("N88").Value = -3919.5

$ws.Range("H91").Value = 2983.077
$ws.Range("I91").Value = 1490
$ws.Range("J91").Value = 3107.5
$ws.Range("K91").Value = 1490
$ws.Range("L91").Value = 3107.5
$ws.Range("M91").Value = -86
$ws.Range("N91").Value = -5915.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1826.8334
$ws.Range("I86").Value = 1571.4286
$ws.Range("J86").Value = 2184.4
$ws.Range("K86").Value = 1571.4286
$ws.Range("L86").Value = 2184.4
$ws.Range("M86").Value = -448.4286
$ws.Range("N86").Value = -4430.4

$ws.Range("H89").Value = 1826.8334
$ws.Range("I89").Value = 1571.4286
$ws.Range("J89").Value = 2184.4
$ws.Range("K89").Value = 7857.143
$ws.Range("L89").Value = 10922
$ws.Range("M89").Value = -2241.143
$ws.Range("N89").Value = -22154

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 22750.092
$ws.Range("I39").Value = 2812.75
$ws.Range("J39").Value = 34142.855
$ws.Range("K39").Value = 2812.75
$ws.Range("L39").Value = 34142.855
$ws.Range("M39").Value = -2421.75
$ws.Range("N39").Value = -34924.855

$ws.Range("H49").Value = 22750.092
$ws.Range("I49").Value = 2812.75
$ws.Range("J49").Value = 34142.855
$ws.Range("K49").Value = 2812.75
$ws.Range("L49").Value = 34142.855
$ws.Range("M49").Value = -2630.75
$ws.Range("N49").Value = -34506.855

$ws.Range("H134").Value = 1315.9615
$ws.Range("I134").Value = 1261.8
$ws.Range("J134").Value = 1389.8182
$ws.Range("K134").Value = 3785.4
$ws.Range("L134").Value = 4169.4546
$ws.Range("M134").Value = -1250.4
$ws.Range("N134").Value = -9239.454600000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1535.75
$ws.Range("I13").Value = 1713.6666
$ws.Range("J13").Value = 1002
$ws.Range("K13").Value = 5140.9998
$ws.Range("L13").Value = 3006
$ws.Range("M13").Value = -4972.9998
$ws.Range("N13").Value = -3342

$ws.Range("H64").Value = 230.5
$ws.Range("I64").Value = 240.66667
$ws.Range("J64").Value = 200
$ws.Range("K64").Value = 722.00001
$ws.Range("L64").Value = 600
$ws.Range("M64").Value = -452.00001
$ws.Range("N64").Value = -1140

$ws.Range("H67").Value = 230.5
$ws.Range("I67").Value = 240.66667
$ws.Range("J67").Value = 200
$ws.Range("K67").Value = 722.00001
$ws.Range("L67").Value = 600
$ws.Range("M67").Value = 213.99999
$ws.Range("N67").Value = -2472

$ws.Range("H76").Value = 6928.5713
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 7250
$ws.Range("K76").Value = 15000
$ws.Range("L76").Value = 21750
$ws.Range("M76").Value = -14617
$ws.Range("N76").Value = -22516

$ws.Range("H79").Value = 6928.5713
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 7250
$ws.Range("K79").Value = 15000
$ws.Range("L79").Value = 21750
$ws.Range("M79").Value = -13674
$ws.Range("N79").Value = -24402

$ws.Range("H99").Value = 100
$ws.Range("I99").Value = 100
$ws.Range("K99").Value = 300
$ws.Range("M99").Value = 1946

$ws.Range("H106").Value = 8644.736999999999
$ws.Range("I106").Value = 500
$ws.Range("J106").Value = 8864.865
$ws.Range("K106").Value = 1500
$ws.Range("L106").Value = 26594.595
$ws.Range("M106").Value = -554
$ws.Range("N106").Value = -28486.595

$ws.Range("H122").Value = 1064242.9
$ws.Range("I122").Value = 300.5
$ws.Range("J122").Value = 5556444
$ws.Range("K122").Value = 2704.5
$ws.Range("L122").Value = 50007996
$ws.Range("M122").Value = -254.5
$ws.Range("N122").Value = -50012896

$ws.Range("H125").Value = 10996.667
$ws.Range("J125").Value = 15995
$ws.Range("L125").Value = 47985
$ws.Range("N125").Value = -57825

$ws.Range("H131").Value = 803.6607
$ws.Range("J131").Value = 975.53656
$ws.Range("L131").Value = 2926.60968
$ws.Range("N131").Value = -13006.60968

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1026.2307
$ws.Range("I46").Value = 861.1111
$ws.Range("J46").Value = 1397.75
$ws.Range("K46").Value = 861.1111
$ws.Range("L46").Value = 1397.75
$ws.Range("M46").Value = -673.1111
$ws.Range("N46").Value = -1773.75

$ws.Range("H82").Value = 1575
$ws.Range("I82").Value = 906.25
$ws.Range("J82").Value = 2110
$ws.Range("K82").Value = 906.25
$ws.Range("L82").Value = 2110
$ws.Range("M82").Value = -545.25
$ws.Range("N82").Value = -2832

$ws.Range("H85").Value = 1575
$ws.Range("I85").Value = 906.25
$ws.Range("J85").Value = 2110
$ws.Range("K85").Value = 906.25
$ws.Range("L85").Value = 2110
$ws.Range("M85").Value = 341.75
$ws.Range("N85").Value = -4606

$ws.Range("H132").Value = 2713.8076
$ws.Range("I132").Value = 1285.7059
$ws.Range("J132").Value = 5411.3335
$ws.Range("K132").Value = 3857.1177
$ws.Range("L132").Value = 16234.0005
$ws.Range("M132").Value = -1327.1177
$ws.Range("N132").Value = -21294.0005

$ws.Range("H136").Value = 1405.5614
$ws.Range("I136").Value = 672.08105
$ws.Range("J136").Value = 2762.5
$ws.Range("K136").Value = 2016.24315
$ws.Range("L136").Value = 8287.5
$ws.Range("M136").Value = 533.75685
$ws.Range("N136").Value = -13387.5

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 569.2
$ws.Range("I113").Value = 534.6667
$ws.Range("J113").Value = 584
$ws.Range("K113").Value = 1604.0001
$ws.Range("L113").Value = 1752
$ws.Range("M113").Value = 565.9999
$ws.Range("N113").Value = -6092

$ws.Range("H132").Value = 1458.9706
$ws.Range("I132").Value = 1339.8572
$ws.Range("J132").Value = 1651.3846
$ws.Range("K132").Value = 4019.5716
$ws.Range("L132").Value = 4954.1538
$ws.Range("M132").Value = -1489.5716
$ws.Range("N132").Value = -10014.1538
